$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "TZP"
$ws.Range("B12").Value = 8163
$ws.Range("C12").Value = 9647.915290584975
$ws.Range("D12").Value = 0.6916445516327471
